# The citation paragraph contains the Duo Security UUID article URL split
# across three separate runs (an artifact of Word's "ignore once" grammar
# check wrapping "#:~" in <w:proofErr> tags):
#
#   run A: "https://duo.com/labs/tech-notes/breaking-down-uuids"
#   proofErr gramStart
#   run B: "#:~"
#   proofErr gramEnd
#   run C: ":text=UUIDs%20are%20generally%20used%20for,physical%20hardware%20within%20an%20organization"
#
# The fix merges A+B+C into a single run (dropping the now-pointless
# proofErr markers) while leaving the surrounding "... from " run and the
# trailing ". " run untouched.

$d = $word.ActiveDocument

# Locate the split URL via Find - Word's Find matches across run / proofErr
# boundaries since it operates on the story's plain text.
$urlRange = $d.Content
$find = $urlRange.Find
$find.ClearFormatting()
$matched = $find.Execute(
    "https://duo.com/labs/tech-notes/breaking-down-uuids#:~:text=UUIDs%20are%20generally%20used%20for,physical%20hardware%20within%20an%20organization",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $matched) {
    throw "Could not find the split duo.com URL citation"
}

# Pull in the trailing ". " run too (one extra char past the match) so the
# replacement range's right edge lands inside a run rather than exactly on
# a run boundary - doing the swap that way keeps the untouched runs on
# either side from being disturbed/reordered by the host's range-splice.
$spliceStart = $urlRange.Start
$spliceEnd = $urlRange.End + 2
$splice = $d.Range($spliceStart, $spliceEnd)

if ($splice.Text -ne ($urlRange.Text + ". ")) {
    throw "Unexpected trailing text after the URL citation: [$($splice.Text)]"
}

$replacementPkg = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidR="009C38F7" w:rsidRPr="00E322E7"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>https://duo.com/labs/tech-notes/breaking-down-uuids#:~:text=UUIDs%20are%20generally%20used%20for,physical%20hardware%20within%20an%20organization</w:t></w:r><w:r w:rsidR="007A6013" w:rsidRPr="007A6013"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p></w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
"@

$splice.InsertXML($replacementPkg)

Write-Output "Merged split duo.com URL citation into a single run."
